# Applies the Alvearie alvearie-fhir-ig gh-pages deployment update for
# StructureDefinition-financial-system.xlsx:
#   - Metadata sheet: Version bump 5.0.0 -> 6.0.0, Date refresh, Publisher set,
#     the duplicated "Contact" row replaced by a new "Jurisdiction" row.
#   - Elements sheet: root Extension row's Short/Definition text updated to
#     describe the financial-system extension specifically.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# The old sheet had two identical "Contact" / "No display for ContactDetail"
# rows (10 and 11). Remove the duplicate (row 10); this shifts every row
# below it up by one, turning the former row 11 "Contact" row into row 10,
# which we then overwrite with the new "Jurisdiction" metadata row.
$meta.Rows.Item(10).Delete()

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refreshed publication timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously blank
$meta.Range("B9").Value = "Alvearie Team"

# New Jurisdiction row replacing the old duplicate Contact row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

$elements = $wb.Worksheets.Item("Elements")

# Root Extension element's Short / Definition text
$elements.Range("K2").Value = "Financial System"
$elements.Range("L2").Value = "Customer-specific code for the financial system"
